$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for several rows as part of
# "repull data, push all data, mean calculation"
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = 3
$ws.Range("F10").Value = -1
